$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the "Total" row (currently row 25) so the new
# task entry becomes row 24 and the Total row shifts down to row 26.
$ws.Rows.Item(24).Insert()

# Fill in the new task row, matching the formatting of the row above it.
$ws.Range("A24").Value = "Added database to the gallery page"

$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value = 43372

$ws.Range("C23").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = 0.5

# Update the Total row's SUM formula to include the new row (now row 25
# is the last data row, Total moved to row 26).
$ws.Range("C26").Formula = "=SUM(C2:C25)"

# Update the active cell selection to match the target state.
$ws.Range("C25").Select()
